# Update the multiplication problems in the table to the new values
# as specified by the commit diff. Each old value is unique in the
# document, so a straightforward Find/Replace for each pair is safe.

$d = $word.ActiveDocument

$replacements = @(
    @("367×8=", "605×3="),
    @("854×8=", "962×5="),
    @("664×8=", "952×9="),
    @("734×8=", "945×2="),
    @("601×5=", "762×9="),
    @("818×7=", "977×8="),
    @("233×3=", "310×8="),
    @("546×2=", "611×8="),
    @("967×6=", "210×7="),
    @("735×4=", "832×6="),
    @("562×2=", "695×2="),
    @("142×5=", "411×9="),
    @("636×6=", "400×2="),
    @("667×3=", "571×4="),
    @("619×5=", "922×2="),
    @("534×2=", "285×8="),
    @("521×9=", "976×3="),
    @("886×9=", "199×8="),
    @("840×7=", "918×7="),
    @("118×2=", "125×7="),
    @("305×8=", "210×2="),
    @("886×6=", "453×3="),
    @("102×5=", "644×2="),
    @("341×5=", "638×8="),
    @("357×6=", "369×6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]

    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
